$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'20.391.56"
$ws.Range("E2").Value = "  +2.29%  "
$ws.Range("D3").Value = "'1.463.17"
$ws.Range("E3").Value = "  +3.89%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").Value = "'0.9460"
$ws.Range("E5").Value = "  -5.60%  "
$ws.Range("D6").Value = "'274.71"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'0.3647"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "'0.3073"
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("D9").Value = "'39.75"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").Value = "'0.06557"
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").Value = "'0.9990"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "'17.97"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "'5.394"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "'6.114"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "'1.460.55"
$ws.Range("E17").Value = "  +3.53%  "
$ws.Range("D18").Value = "'0.9625"
$ws.Range("E18").Value = "  -3.91%  "
$ws.Range("D19").Value = "'0.05753"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("D20").Value = "'69.58"
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("D21").Value = "'5.413"
$ws.Range("E21").Value = "  -3.51%  "
$ws.Range("D22").Value = "'14.40"
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("D24").Value = "'2.224"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("D25").Value = "'20.415.68"
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("D26").Value = "'141.01"
$ws.Range("E26").Value = "  +6.28%  "
$ws.Range("E27").Value = "  -7.74%  "
$ws.Range("D28").Value = "'17.10"
$ws.Range("E28").Value = "  -1.12%  "
$ws.Range("D29").Value = "'1.610.40"
$ws.Range("E29").Value = "  +2.61%  "
$ws.Range("D30").Value = "'111.89"
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("D31").Value = "'3.812"
$ws.Range("E31").Value = "  -4.24%  "
$ws.Range("D32").Value = "'4.861"
$ws.Range("E32").Value = "  -7.81%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "'0.07800"
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7865"
$ws.Range("E34").Value = "  -3.48%  "
$ws.Range("D35").Value = "'1.502"
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("D36").Value = "'0.05690"
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("D37").Value = "'4.652"
$ws.Range("E37").Value = "  -5.21%  "
$ws.Range("E38").Value = "  +3.01%  "
$ws.Range("D39").Value = "'0.02029"
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D40").Value = "'0.9520"
$ws.Range("E40").Value = "  -4.89%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'10.33"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").Value = "'0.1858"
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("D43").Value = "'7.409"
$ws.Range("E43").Value = "  -10.94%  "
$ws.Range("D44").Value = "'0.5254"
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("D45").Value = "'3.485"
$ws.Range("D46").Value = "'11.83"
$ws.Range("E46").Value = "  -4.23%  "
$ws.Range("D47").Value = "'117.06"
$ws.Range("E47").Value = "  +2.04%  "
$ws.Range("D48").Value = "'0.5136"
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("D49").Value = "'1.746"
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("E50").Value = "  +3.85%  "
$ws.Range("D51").Value = "'0.9874"
$ws.Range("E51").Value = "  -1.52%  "
